$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "Arab News"
$ws.Range("B2").Value = "2001-12-19T00:00:00UTC"
$ws.Range("E2").Value = "http://www.arabnews.com/?page=1&section=0&article=61839&d=9&m=4&y=2005&pix=kingdom.jpg&category=Kingdom"

$ws.Range("A3").Value = "Djibouti incumbent wins one-man poll"
$ws.Range("B3").Value = "1-01-01T00:00:00UTC"
$ws.Range("E3").Value = "http://www.afrol.com/articles/16113"
